# feat: add 2022-Q4 data
#
# 1. The sheet currently named "2021-Q4" becomes "2022-Q4" and gains two new
#    fund rows at the top of its data (old rows shift down).
# 2. A brand-new sheet named "2021-Q4" is appended, holding an exact copy of
#    what the "2021-Q4" sheet used to contain before today's edit.
# 3. The summary sheet ("总计") gets its existing 2021-Q4 row updated to
#    describe 2022-Q4, and a new row is appended with the original 2021-Q4
#    figures.

$wb = $excel.ActiveWorkbook

$wsSummary = $wb.Worksheets.Item(1)
$wsQ = $wb.Worksheets.Item(2)

# --- Step 1: duplicate the existing quarter sheet so the old data survives
# under its own "2021-Q4" tab, placed right after the sheet being updated.
$wsQ.Copy($null, $wsQ)
$wsOld = $wb.Worksheets.Item(3)
$wsOld.Name = "2021-Q4-temp"

# Rename the original (soon to hold the new quarter's data).
$wsQ.Name = "2022-Q4"
$wsOld.Name = "2021-Q4"

# --- Step 2: insert two new fund rows at the top of the "2022-Q4" sheet data
$wsQ.Rows.Item(2).Insert()
$wsQ.Rows.Item(2).Insert()

# Re-apply the bold/centered/bordered look used by the sibling data rows to
# the freshly inserted A-column cells (row insert only carried it to B:H).
$srcFmt = $wsQ.Range("A4")
$fmtTargets = @($wsQ.Range("A2"), $wsQ.Range("A3"))
foreach ($cell in $fmtTargets) {
    $cell.Font.Bold = $srcFmt.Font.Bold
    $cell.HorizontalAlignment = $srcFmt.HorizontalAlignment
    $cell.VerticalAlignment = $srcFmt.VerticalAlignment
    $cell.Borders.Item(1).Weight = $srcFmt.Borders.Item(1).Weight
    $cell.Borders.Item(2).Weight = $srcFmt.Borders.Item(2).Weight
    $cell.Borders.Item(3).Weight = $srcFmt.Borders.Item(3).Weight
    $cell.Borders.Item(4).Weight = $srcFmt.Borders.Item(4).Weight
}

# Columns B:G hold text values (fund code/name/size/position/weight/value);
# force text storage so things like "013360" keep their leading zero and
# "0.50" keeps its trailing zero instead of becoming the number 0.5.
$wsQ.Range("B2:G3").NumberFormat = "@"

$wsQ.Range("A2").Value = 0
$wsQ.Range("B2").Value = "160323"
$wsQ.Range("C2").Value = "华夏磐泰混合（LOF）A"
$wsQ.Range("D2").Value = "4.95"
$wsQ.Range("E2").Value = "28.84"
$wsQ.Range("F2").Value = "0.50"
$wsQ.Range("G2").Value = "0.0248"
$wsQ.Range("H2").Value = 9

$wsQ.Range("A3").Value = 1
$wsQ.Range("B3").Value = "013360"
$wsQ.Range("C3").Value = "华夏磐泰混合（LOF）C"
$wsQ.Range("D3").Value = "3.48"
$wsQ.Range("E3").Value = "28.84"
$wsQ.Range("F3").Value = "0.50"
$wsQ.Range("G3").Value = "0.0174"
$wsQ.Range("H3").Value = 9

# --- Step 3: update the summary sheet
$wsSummary.Range("B2").Value = "2022-Q4"
$wsSummary.Range("D2").Value = 0.04

$wsSummary.Range("A3").Value = 1
$wsSummary.Range("B3").Value = "2021-Q4"
$wsSummary.Range("C3").Value = 2
$wsSummary.Range("D3").Value = 0.09

# Match the bold/centered/bordered style of the row above for the new A3 cell.
$srcFmt2 = $wsSummary.Range("A2")
$wsSummary.Range("A3").Font.Bold = $srcFmt2.Font.Bold
$wsSummary.Range("A3").HorizontalAlignment = $srcFmt2.HorizontalAlignment
$wsSummary.Range("A3").VerticalAlignment = $srcFmt2.VerticalAlignment
$wsSummary.Range("A3").Borders.Item(1).Weight = $srcFmt2.Borders.Item(1).Weight
$wsSummary.Range("A3").Borders.Item(2).Weight = $srcFmt2.Borders.Item(2).Weight
$wsSummary.Range("A3").Borders.Item(3).Weight = $srcFmt2.Borders.Item(3).Weight
$wsSummary.Range("A3").Borders.Item(4).Weight = $srcFmt2.Borders.Item(4).Weight

# Leave the workbook focused the way it started: summary sheet, cell A1.
$wsSummary.Activate() | Out-Null
$wsSummary.Range("A1").Select() | Out-Null
